$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F "time_taken" to the panel export, recording when each row was processed.

# Header cell F1: set text, then copy the formatting of the other header cells (E1) onto it
# so it keeps the bold font / border / centered-top alignment used by the rest of row 1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data rows F2:F110 - per-row "time_taken" timestamp values (plain text, not dates).
$timeTaken = @(
    "2021-10-05 10:50:17.643697",
    "2021-10-05 10:50:17.643708",
    "2021-10-05 10:50:17.643712",
    "2021-10-05 10:50:17.643714",
    "2021-10-05 10:50:17.643717",
    "2021-10-05 10:50:17.643720",
    "2021-10-05 10:50:17.643722",
    "2021-10-05 10:50:17.643724",
    "2021-10-05 10:50:17.643727",
    "2021-10-05 10:50:17.643730",
    "2021-10-05 10:50:17.643732",
    "2021-10-05 10:50:17.643735",
    "2021-10-05 10:50:17.643738",
    "2021-10-05 10:50:17.643740",
    "2021-10-05 10:50:17.643742",
    "2021-10-05 10:50:17.643745",
    "2021-10-05 10:50:17.643748",
    "2021-10-05 10:50:17.643750",
    "2021-10-05 10:50:17.643753",
    "2021-10-05 10:50:17.643755",
    "2021-10-05 10:50:17.643758",
    "2021-10-05 10:50:17.643760",
    "2021-10-05 10:50:17.643763",
    "2021-10-05 10:50:17.643765",
    "2021-10-05 10:50:17.643768",
    "2021-10-05 10:50:17.643771",
    "2021-10-05 10:50:17.643773",
    "2021-10-05 10:50:17.643776",
    "2021-10-05 10:50:17.643778",
    "2021-10-05 10:50:17.643780",
    "2021-10-05 10:50:17.643783",
    "2021-10-05 10:50:17.643785",
    "2021-10-05 10:50:17.643789",
    "2021-10-05 10:50:17.643791",
    "2021-10-05 10:50:17.643794",
    "2021-10-05 10:50:17.643796",
    "2021-10-05 10:50:17.643799",
    "2021-10-05 10:50:17.643801",
    "2021-10-05 10:50:17.643804",
    "2021-10-05 10:50:17.643806",
    "2021-10-05 10:50:17.643809",
    "2021-10-05 10:50:17.643811",
    "2021-10-05 10:50:17.643814",
    "2021-10-05 10:50:17.643816",
    "2021-10-05 10:50:17.643819",
    "2021-10-05 10:50:17.643821",
    "2021-10-05 10:50:17.643824",
    "2021-10-05 10:50:17.643826",
    "2021-10-05 10:50:17.643829",
    "2021-10-05 10:50:17.643831",
    "2021-10-05 10:50:17.643834",
    "2021-10-05 10:50:17.643836",
    "2021-10-05 10:50:17.643839",
    "2021-10-05 10:50:17.643842",
    "2021-10-05 10:50:17.643844",
    "2021-10-05 10:50:17.643847",
    "2021-10-05 10:50:17.643849",
    "2021-10-05 10:50:17.643852",
    "2021-10-05 10:50:17.643854",
    "2021-10-05 10:50:17.643857",
    "2021-10-05 10:50:17.643859",
    "2021-10-05 10:50:17.643862",
    "2021-10-05 10:50:17.643864",
    "2021-10-05 10:50:17.643866",
    "2021-10-05 10:50:17.643870",
    "2021-10-05 10:50:17.643873",
    "2021-10-05 10:50:17.643875",
    "2021-10-05 10:50:17.643878",
    "2021-10-05 10:50:17.643880",
    "2021-10-05 10:50:17.643883",
    "2021-10-05 10:50:17.643885",
    "2021-10-05 10:50:17.643888",
    "2021-10-05 10:50:17.643890",
    "2021-10-05 10:50:17.643893",
    "2021-10-05 10:50:17.643895",
    "2021-10-05 10:50:17.643898",
    "2021-10-05 10:50:17.643902",
    "2021-10-05 10:50:17.643905",
    "2021-10-05 10:50:17.643907",
    "2021-10-05 10:50:17.643910",
    "2021-10-05 10:50:17.643912",
    "2021-10-05 10:50:17.643915",
    "2021-10-05 10:50:17.643917",
    "2021-10-05 10:50:17.643920",
    "2021-10-05 10:50:17.643922",
    "2021-10-05 10:50:17.643925",
    "2021-10-05 10:50:17.643927",
    "2021-10-05 10:50:17.643929",
    "2021-10-05 10:50:17.643932",
    "2021-10-05 10:50:17.643934",
    "2021-10-05 10:50:17.643937",
    "2021-10-05 10:50:17.643940",
    "2021-10-05 10:50:17.643943",
    "2021-10-05 10:50:17.643946",
    "2021-10-05 10:50:17.643948",
    "2021-10-05 10:50:17.643951",
    "2021-10-05 10:50:17.643953",
    "2021-10-05 10:50:17.643956",
    "2021-10-05 10:50:17.643958",
    "2021-10-05 10:50:17.643961",
    "2021-10-05 10:50:17.643963",
    "2021-10-05 10:50:17.643966",
    "2021-10-05 10:50:17.643968",
    "2021-10-05 10:50:17.643971",
    "2021-10-05 10:50:17.643973",
    "2021-10-05 10:50:17.643976",
    "2021-10-05 10:50:17.643978",
    "2021-10-05 10:50:17.643981",
    "2021-10-05 10:50:17.643986"
)

$firstDataRow = 2
for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $firstDataRow + $i
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
